$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - first data sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 579
$ws1.Range("F5").Value = 291
$ws1.Range("F6").Value = 1112
$ws1.Range("F7").Value = 1444
$ws1.Range("F9").Value = 112
$ws1.Range("F10").Value = 753
$ws1.Range("F12").Value = 173
$ws1.Range("F14").Value = 446
$ws1.Range("F15").Value = 1386
$ws1.Range("F16").Value = 121
$ws1.Range("F17").Value = 119
$ws1.Range("F18").Value = 279
$ws1.Range("F19").Value = 5210
$ws1.Range("F20").Value = 62
$ws1.Range("F21").Value = 661
$ws1.Range("F22").Value = 1008
$ws1.Range("F23").Value = 37
$ws1.Range("F24").Value = 238
$ws1.Range("F26").Value = 5956
$ws1.Range("F27").Value = 72
$ws1.Range("F28").Value = 123
$ws1.Range("F29").Value = 116
$ws1.Range("F31").Value = 14646
$ws1.Range("F32").Value = 1453
$ws1.Range("F33").Value = 223
$ws1.Range("F34").Value = 106
$ws1.Range("F35").Value = 89
$ws1.Range("F36").Value = 9426
$ws1.Range("F37").Value = 638
$ws1.Range("F38").Value = 4221
$ws1.Range("F39").Value = 155

# Sheet "全部类型" (All Types) - combined data sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 579
$ws4.Range("F5").Value = 291
$ws4.Range("F6").Value = 1112
$ws4.Range("F7").Value = 1444
$ws4.Range("F9").Value = 112
$ws4.Range("F10").Value = 753
$ws4.Range("F12").Value = 173
$ws4.Range("F14").Value = 446
$ws4.Range("F15").Value = 1386
$ws4.Range("F16").Value = 121
$ws4.Range("F17").Value = 119
$ws4.Range("F18").Value = 279
$ws4.Range("F20").Value = 5210
$ws4.Range("F21").Value = 62
$ws4.Range("F22").Value = 661
$ws4.Range("F24").Value = 1008
$ws4.Range("F25").Value = 37
$ws4.Range("F26").Value = 238
$ws4.Range("F29").Value = 5956
$ws4.Range("F30").Value = 72
$ws4.Range("F31").Value = 123
$ws4.Range("F32").Value = 116
$ws4.Range("F34").Value = 14646
$ws4.Range("F35").Value = 1453
$ws4.Range("F36").Value = 223
$ws4.Range("F37").Value = 106
$ws4.Range("F38").Value = 89
$ws4.Range("F39").Value = 9428
$ws4.Range("F40").Value = 638
$ws4.Range("F41").Value = 4221
$ws4.Range("F42").Value = 155
